# Update column F (dSF) values for the specified rows on the active sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    3  = 0
    4  = -5
    5  = 0
    6  = -4
    8  = 6
    9  = 4
    10 = 2
    12 = 4
    14 = 2
    15 = 1
    16 = 4
    17 = 5
    18 = 1
    19 = 6
    20 = 5
    21 = 4
    22 = 1
    23 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
